$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.043.60'
$ws.Range('E2').Value = '  -3.74%  '
$ws.Range('D3').Value = '1.745.73'
$ws.Range('E3').Value = '  -4.31%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.97'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5805'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.98%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.002'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2710'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.31%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '23.15'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06582'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07502'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.13%  '
$ws.Range('D12').Value = '1.730.66'
$ws.Range('E12').Value = '  -5.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.714'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.30%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6043'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.28%  '
$ws.Range('D15').Value = '1.982.10'
$ws.Range('E15').Value = '  -4.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '74.12'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.81%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008649'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -11.79%  '
$ws.Range('D18').Value = '28.020.64'
$ws.Range('E18').Value = '  -2.75%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.320'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.22%  '
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '204.74'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.25'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.616'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.003'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '150.50'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.64%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.017'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1236'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.10'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.392'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.86%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.06120'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.387'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.75%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.742'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.716'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.93%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.677'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.73%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.036'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6352'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.457'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.710'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01674'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.14%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.278'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.68%  '
$ws.Range('D41').Value = '1.125.42'
$ws.Range('E41').Value = '  -1.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8656'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.25%  '
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.47'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.54%  '
$ws.Range('D45').Value = '1.894.36'
$ws.Range('E45').Value = '  -4.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '59.36'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.573'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.56%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.00000000106'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.63%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.238'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.71%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05380'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.06%  '
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.287'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.70%  '
